function Set-RowValues {
    param($ws, $row, $vals)
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $vals[$i]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-RowValues $ws 2 @("ECs", "Gnai2", "Fpr1", "ECs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 2223.956981120573, 20015.61283008516, 0.2527940559544053, 0.2527940559544053)
Set-RowValues $ws 3 @("ECs", "Gnai2", "Fpr1", "FAPs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 381.137995715943, 3430.241961443487, 0.04332341885804775, 0.04332341885804775)
Set-RowValues $ws 4 @("ECs", "Gnai2", "Fpr1", "M2", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 618.9991033874379, 5570.991930486941, 0.07036075576363278, 0.07036075576363277)
Set-RowValues $ws 5 @("ECs", "Gnai2", "Fpr1", "sCs", 3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 216.057378056275, 1944.516402506475, 0.02455893768691514, 0.02455893768691514)
Set-RowValues $ws 6 @("FAPs", "Gnai2", "Fpr1", "ECs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 895.8499668686599, 8062.649701817939, 0.1018300032661784, 0.1018300032661784)
Set-RowValues $ws 7 @("FAPs", "Gnai2", "Fpr1", "FAPs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 153.5292560661287, 1381.763304595158, 0.01745145417743806, 0.01745145417743806)
Set-RowValues $ws 8 @("FAPs", "Gnai2", "Fpr1", "M2", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 249.3439985435142, 2244.095986891628, 0.02834258092885608, 0.02834258092885608)
Set-RowValues $ws 9 @("FAPs", "Gnai2", "Fpr1", "sCs", 3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 87.03180709723888, 783.2862638751499, 0.00989278286402232, 0.00989278286402232)
Set-RowValues $ws 10 @("M2", "Gnai2", "Fpr1", "ECs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 1729.694731315046, 15567.25258183541, 0.1966119625532404, 0.1966119625532405)
Set-RowValues $ws 11 @("M2", "Gnai2", "Fpr1", "FAPs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 296.4321651409233, 2667.88948626831, 0.03369502646744648, 0.03369502646744649)
Set-RowValues $ws 12 @("M2", "Gnai2", "Fpr1", "M2", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 481.4299453213844, 4332.86950789246, 0.05472346343425164, 0.05472346343425164)
Set-RowValues $ws 13 @("M2", "Gnai2", "Fpr1", "sCs", 3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 168.0398099685278, 1512.35828971675, 0.01910084839066795, 0.01910084839066795)
Set-RowValues $ws 14 @("sCs", "Gnai2", "Fpr1", "ECs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 837.8272030100496, 7540.444827090448, 0.09523463746638121, 0.09523463746638124)
Set-RowValues $ws 15 @("sCs", "Gnai2", "Fpr1", "FAPs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 143.585412677653, 1292.268714098877, 0.01632115151273382, 0.01632115151273382)
Set-RowValues $ws 16 @("sCs", "Gnai2", "Fpr1", "M2", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 233.1943881376313, 2098.749493238682, 0.02650687747269949, 0.02650687747269949)
Set-RowValues $ws 17 @("sCs", "Gnai2", "Fpr1", "sCs", 3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 81.39489670135832, 732.5540703122249, 0.00925204320308328, 0.009252043203083283)

Write-Output "done"
